# Apply the edits described by the diff:
#  - C73: 0.15 -> 0.25
#  - C74: 1.75 -> 2.25
#  - D75: new note "Start Lesson 4"
#  - Shared string used by D74 ("further exploration and 5 small problems")
#    changes text to "further exploration and 9 small problems"
#  - sheetView: topLeftCell B67 -> A64, selection C75 -> C74

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the hours logged on 2021-11-07 and 2021-11-08
$ws.Range("C73").Value = 0.25
$ws.Range("C74").Value = 2.25

# Rename/expand the milestone text referenced from D74 (shared string swap)
$ws.Range("D74").Value = "further exploration and 9 small problems"

# Add the new milestone note for row 75
$ws.Range("D75").Value = "Start Lesson 4"

# Update the active selection / scroll position to match the saved view
$ws.Range("C74").Select()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
